$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 865.44446
$ws.Range("J70").Value = 878
$ws.Range("L70").Value = 2634
$ws.Range("N70").Value = -3174
$ws.Range("H73").Value = 865.44446
$ws.Range("J73").Value = 878
$ws.Range("L73").Value = 2634
$ws.Range("N73").Value = -4506
$ws.Range("H88").Value = 1746.0625
$ws.Range("J88").Value = 1828.3334
$ws.Range("L88").Value = 1828.3334
$ws.Range("N88").Value = -2640.3334
$ws.Range("H91").Value = 1746.0625
$ws.Range("J91").Value = 1828.3334
$ws.Range("L91").Value = 1828.3334
$ws.Range("N91").Value = -4636.3334
$ws.Range("H116").Value = 3505.25
$ws.Range("I116").Value = 1466
$ws.Range("K116").Value = 1466
$ws.Range("M116").Value = 1976

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1294.9
$ws.Range("I2").Value = 1081.4
$ws.Range("J2").Value = 2362.4
$ws.Range("K2").Value = 1081.4
$ws.Range("L2").Value = 2362.4
$ws.Range("M2").Value = -968.4000000000001
$ws.Range("N2").Value = -2588.4
$ws.Range("H45").Value = 1579.2
$ws.Range("I45").Value = 1517.8077
$ws.Range("J45").Value = 1693.2142
$ws.Range("K45").Value = 1517.8077
$ws.Range("L45").Value = 1693.2142
$ws.Range("M45").Value = -1140.8077
$ws.Range("N45").Value = -2447.2142
$ws.Range("H63").Value = 1798.4286
$ws.Range("I63").Value = 1798.4286
$ws.Range("K63").Value = 1798.4286
$ws.Range("M63").Value = -1112.4286
$ws.Range("H66").Value = 1798.4286
$ws.Range("I66").Value = 1798.4286
$ws.Range("K66").Value = 8992.143
$ws.Range("M66").Value = -5560.143
$ws.Range("H74").Value = 76923960
$ws.Range("I74").Value = 111111610
$ws.Range("K74").Value = 111111610
$ws.Range("M74").Value = -111110736
$ws.Range("H77").Value = 76923960
$ws.Range("I77").Value = 111111610
$ws.Range("K77").Value = 555558050
$ws.Range("M77").Value = -555553682
$ws.Range("H80").Value = 47971
$ws.Range("J80").Value = 47971
$ws.Range("L80").Value = 47971
$ws.Range("N80").Value = -49967
$ws.Range("H83").Value = 47971
$ws.Range("J83").Value = 47971
$ws.Range("L83").Value = 143913
$ws.Range("N83").Value = -153897
$ws.Range("H113").Value = 35000
$ws.Range("J113").Value = 35000
$ws.Range("L113").Value = 35000
$ws.Range("N113").Value = -43678
$ws.Range("H116").Value = 1294.9
$ws.Range("I116").Value = 1081.4
$ws.Range("J116").Value = 2362.4
$ws.Range("K116").Value = 1081.4
$ws.Range("L116").Value = 2362.4
$ws.Range("M116").Value = 1212.6
$ws.Range("N116").Value = -6950.4
$ws.Range("H124").Value = 15224.833
$ws.Range("J124").Value = 15224.833
$ws.Range("L124").Value = 15224.833
$ws.Range("N124").Value = -25044.833

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1294.9
$ws.Range("I3").Value = 1081.4
$ws.Range("J3").Value = 2362.4
$ws.Range("K3").Value = 1081.4
$ws.Range("L3").Value = 2362.4
$ws.Range("M3").Value = -967.4000000000001
$ws.Range("N3").Value = -2590.4
$ws.Range("H50").Value = 23590
$ws.Range("J50").Value = 23590
$ws.Range("L50").Value = 23590
$ws.Range("N50").Value = -24738
$ws.Range("H134").Value = 41382.594
$ws.Range("I134").Value = 48310
$ws.Range("J134").Value = 1550
$ws.Range("K134").Value = 144930
$ws.Range("L134").Value = 4650
$ws.Range("M134").Value = -142395
$ws.Range("N134").Value = -9720

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 303.0909
$ws.Range("I22").Value = 171.66667
$ws.Range("J22").Value = 460.8
$ws.Range("K22").Value = 171.66667
$ws.Range("L22").Value = 460.8
$ws.Range("M22").Value = 178.33333
$ws.Range("N22").Value = -1160.8
$ws.Range("H31").Value = 12347.19
$ws.Range("J31").Value = 4481.421
$ws.Range("L31").Value = 4481.421
$ws.Range("N31").Value = -5071.421
$ws.Range("H34").Value = 12347.19
$ws.Range("J34").Value = 4481.421
$ws.Range("L34").Value = 4481.421
$ws.Range("N34").Value = -4885.421
$ws.Range("H62").Value = 125003850
$ws.Range("I62").Value = 166670800
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 166670800
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -166670176
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 125003850
$ws.Range("I65").Value = 166670800
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 833354000
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -833350880
$ws.Range("N65").Value = -21240
$ws.Range("H99").Value = 19235080
$ws.Range("I99").Value = 3573.1428
$ws.Range("J99").Value = 41671836
$ws.Range("K99").Value = 3573.1428
$ws.Range("L99").Value = 41671836
$ws.Range("M99").Value = -2075.1428
$ws.Range("N99").Value = -41674832
$ws.Range("H121").Value = 7747.857
$ws.Range("I121").Value = 7747.857
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 7747.857
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = -6437.857
$ws.Range("N121").ClearContents() | Out-Null
$ws.Range("H126").Value = 19235080
$ws.Range("I126").Value = 3573.1428
$ws.Range("J126").Value = 41671836
$ws.Range("K126").Value = 10719.4284
$ws.Range("L126").Value = 125015508
$ws.Range("M126").Value = -8249.428400000001
$ws.Range("N126").Value = -125020448

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 785.4400000000001
$ws.Range("J131").Value = 786.30304
$ws.Range("L131").Value = 2358.90912
$ws.Range("N131").Value = -12438.90912

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4668.9653
$ws.Range("I126").Value = 3545
$ws.Range("J126").Value = 7166.6665
$ws.Range("K126").Value = 10635
$ws.Range("L126").Value = 21499.9995
$ws.Range("M126").Value = -8165
$ws.Range("N126").Value = -26439.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 29332.666
$ws.Range("J36").Value = 29332.666
$ws.Range("L36").Value = 29332.666
$ws.Range("N36").Value = -30456.666
$ws.Range("H40").Value = 4143.0586
$ws.Range("I40").Value = 2649.2856
$ws.Range("J40").Value = 5188.7
$ws.Range("K40").Value = 2649.2856
$ws.Range("L40").Value = 5188.7
$ws.Range("M40").Value = -2513.2856
$ws.Range("N40").Value = -5460.7
$ws.Range("H68").Value = 2784.625
$ws.Range("I68").Value = 3133.3333
$ws.Range("J68").Value = 2575.4
$ws.Range("K68").Value = 3133.3333
$ws.Range("L68").Value = 2575.4
$ws.Range("M68").Value = -2384.3333
$ws.Range("N68").Value = -4073.4
$ws.Range("H71").Value = 2784.625
$ws.Range("I71").Value = 3133.3333
$ws.Range("J71").Value = 2575.4
$ws.Range("K71").Value = 15666.6665
$ws.Range("L71").Value = 12877
$ws.Range("M71").Value = -11922.6665
$ws.Range("N71").Value = -20365
$ws.Range("H93").Value = 1253
$ws.Range("I93").Value = 1335.0454
$ws.Range("J93").Value = 892
$ws.Range("K93").Value = 1335.0454
$ws.Range("L93").Value = 892
$ws.Range("M93").Value = -87.04539999999997
$ws.Range("N93").Value = -3388
$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840
$ws.Range("H132").Value = 2826.1482
$ws.Range("I132").Value = 1618.4546
$ws.Range("K132").Value = 4855.3638
$ws.Range("M132").Value = -2325.3638
$ws.Range("H136").Value = 57000.777
$ws.Range("I136").Value = 126500.5
$ws.Range("J136").Value = 1401
$ws.Range("K136").Value = 379501.5
$ws.Range("L136").Value = 4203
$ws.Range("M136").Value = -376951.5
$ws.Range("N136").Value = -9303

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4628
$ws.Range("I62").Value = 3439
$ws.Range("J62").Value = 5371.125
$ws.Range("K62").Value = 3439
$ws.Range("L62").Value = 5371.125
$ws.Range("M62").Value = -2815
$ws.Range("N62").Value = -6619.125
$ws.Range("H65").Value = 4628
$ws.Range("I65").Value = 3439
$ws.Range("J65").Value = 5371.125
$ws.Range("K65").Value = 17195
$ws.Range("L65").Value = 26855.625
$ws.Range("M65").Value = -14075
$ws.Range("N65").Value = -33095.625
$ws.Range("H126").Value = 996.125
$ws.Range("I126").Value = 821.4286
$ws.Range("J126").Value = 1132
$ws.Range("K126").Value = 2464.2858
$ws.Range("L126").Value = 3396
$ws.Range("M126").Value = 5.714200000000346
$ws.Range("N126").Value = -8336
$ws.Range("H132").Value = 3174.5
$ws.Range("I132").Value = 1350
$ws.Range("K132").Value = 4050
$ws.Range("M132").Value = -1520
$ws.Range("H135").Value = 50853.25
$ws.Range("J135").Value = 50853.25
$ws.Range("L135").Value = 50853.25
$ws.Range("N135").Value = -60993.25
$ws.Range("H136").Value = 41669140
$ws.Range("I136").Value = 83335980
$ws.Range("K136").Value = 250007940
$ws.Range("M136").Value = -250005390
